$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (column C) date values for rows 2-10 from 45186 to 45188
# (serial date 2023-09-17 -> 2023-09-19), keeping existing number formatting.
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 3).Value = 45188
}
